# Apply an automatic update that re-orders the species observation rows
# (rows 2,3,5,6,7,8,9 ; row 4 is unchanged) by shuffling the columns
# A,B,D,E,F,G,H,Q,R,S between rows while leaving all other columns
# (C,P,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AT,AW,AX,AY, etc.) untouched per-row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move between rows.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "S")

# Mapping of destination row -> source row (data that ends up in the
# destination row is taken from the source row's original contents).
$rowMap = @{
    2 = 3
    3 = 6
    5 = 7
    6 = 5
    7 = 2
    8 = 9
    9 = 8
}

# Snapshot the original values for the columns that will be rearranged,
# before any writes happen (so later writes don't clobber values we
# still need to read for other rows).
$original = @{}
foreach ($r in 2..9) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $original[$r] = $rowVals
}

# Write the rearranged values into each destination row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcVals[$col]
    }
}
